$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "yQjMV407"
$ws.Range("B2").Value = 231102303
$ws.Range("C2").Value = "mtbhnpn84"
$ws.Range("D2").Value = 'hP!2$z9U'
$ws.Range("F2").Value = "NRFWWfRD"
$ws.Range("G2").Value = "lxwa"
